$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 ("Experimental" row) currently empty -> "false".
# A plain Value="false" assignment is interpreted as a Boolean by the
# engine (mirrors Excel's literal-entry parser), so enter it with a
# leading apostrophe to force text, then re-apply the plain body format
# (copied from the neighboring data cell) so the cell's style/format
# matches the rest of the column instead of picking up a quote-prefix style.
$ws.Range("B7").Value = "'false"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# B8 ("Date" row) -> updated timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
